# Apply updated Betfair back/lay odds per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("I2").Value = 4.2
$ws.Range("P2").Value = 2.22
$ws.Range("AN2").Value = 13.5
# Row 3
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 2.72
$ws.Range("AM3").Value = 110
# Row 4
$ws.Range("H4").Value = 17
$ws.Range("J4").Value = 1.09
# Row 5
$ws.Range("Q5").Value = 2.6
# Row 9
$ws.Range("J9").Value = 5.6
$ws.Range("K9").Value = 5.9
$ws.Range("P9").Value = 2.42
$ws.Range("R9").Value = 1.58
$ws.Range("AC9").Value = 13
$ws.Range("AD9").Value = 48
$ws.Range("AE9").Value = 250
$ws.Range("AF9").Value = 8.199999999999999
$ws.Range("AH9").Value = 50
$ws.Range("AN9").Value = 4.9
# Row 10
$ws.Range("G10").Value = 1.6
$ws.Range("H10").Value = 6
$ws.Range("I10").Value = 6.4
$ws.Range("J10").Value = 4.7
$ws.Range("K10").Value = 5
$ws.Range("N10").Value = 7.2
$ws.Range("Q10").Value = 1.44
$ws.Range("S10").Value = 2.08
$ws.Range("X10").Value = 34
$ws.Range("Y10").Value = 36
$ws.Range("AA10").Value = 170
$ws.Range("AD10").Value = 25
$ws.Range("AE10").Value = 65
$ws.Range("AF10").Value = 13.5
$ws.Range("AM10").Value = 55
$ws.Range("AO10").Value = 1000
# Row 14
$ws.Range("G14").Value = 2.34
$ws.Range("N14").Value = 4.9
$ws.Range("X14").Value = 20
$ws.Range("Y14").Value = 16.5
$ws.Range("AA14").Value = 65
$ws.Range("AF14").Value = 17
$ws.Range("AI14").Value = 1000
$ws.Range("AM14").Value = 1000
# Row 15
$ws.Range("F15").Value = 3.1
$ws.Range("G15").Value = 3.15
$ws.Range("H15").Value = 2.42
$ws.Range("I15").Value = 2.46
$ws.Range("K15").Value = 3.75
$ws.Range("M15").Value = 1.05
$ws.Range("R15").Value = 1.56
$ws.Range("Z15").Value = 21
$ws.Range("AA15").Value = 36
$ws.Range("AH15").Value = 15
$ws.Range("AI15").Value = 34
$ws.Range("AM15").Value = 60
# Row 16
$ws.Range("F16").Value = 1.51
$ws.Range("G16").Value = 1.52
$ws.Range("H16").Value = 6.6
$ws.Range("I16").Value = 7.2
$ws.Range("Q16").Value = 1.57
$ws.Range("R16").Value = 1.66
$ws.Range("U16").Value = 2.3
$ws.Range("X16").Value = 28
$ws.Range("Y16").Value = 32
$ws.Range("AB16").Value = 12
$ws.Range("AD16").Value = 27
$ws.Range("AF16").Value = 12
$ws.Range("AI16").Value = 80
$ws.Range("AM16").Value = 1000
# Row 17
$ws.Range("F17").Value = 1.23
$ws.Range("G17").Value = 1.24
$ws.Range("I17").Value = 17
$ws.Range("J17").Value = 7.6
$ws.Range("Q17").Value = 1.31
$ws.Range("Z17").Value = 200
$ws.Range("AE17").Value = 230
# Row 18
$ws.Range("H18").Value = 12.5
$ws.Range("K18").Value = 7.8
$ws.Range("U18").Value = 2.02
$ws.Range("X18").Value = 60
$ws.Range("AB18").Value = 13
$ws.Range("AC18").Value = 18
# Row 19
$ws.Range("G19").Value = 5.7
$ws.Range("R19").Value = 1.91
$ws.Range("X19").Value = 38
$ws.Range("AN19").Value = 36
# Row 20
$ws.Range("F20").Value = 1.54
$ws.Range("G20").Value = 1.55
$ws.Range("H20").Value = 6.6
$ws.Range("I20").Value = 7
$ws.Range("X20").Value = 24
$ws.Range("Y20").Value = 28
$ws.Range("AD20").Value = 26
$ws.Range("AH20").Value = 22
# Row 21
$ws.Range("F21").Value = 1.83
$ws.Range("H21").Value = 1.33
$ws.Range("I21").Value = 4.9
$ws.Range("K21").Value = 4
# Row 22
$ws.Range("F22").Value = 1.98
$ws.Range("G22").Value = 2.06
$ws.Range("K22").Value = 3.25
$ws.Range("O22").Value = 1.66
$ws.Range("Q22").Value = 3
$ws.Range("T22").Value = 2.5
$ws.Range("U22").Value = 1.56
$ws.Range("W22").Value = 1.94
$ws.Range("AN22").Value = 36
